{"js": "// Apply the three textual insertions described by the diff.\n// Each change appends additional sentence(s)/phrase(s) immediately after an\n// existing, uniquely-identifiable piece of text within the same paragraph.\n\nconst body = context.document.body;\n\n// --- Change 1 --------------------------------------------------------\n// \"...where I had to revise multiple times.\" ->\n// \"...where I had to revise multiple times, and I am proud of the revised\n//  results I was able to get at the end.\"\nlet results1 = body.search(\"where I had to revise multiple times\", { matchCase: true, matchWholeWord: false });\nresults1.load(\"items\");\nawait context.sync();\nif (results1.items.length > 0) {\n  results1.items[0].insertText(\n    \", and I am proud of the revised results I was able to get at the end\",\n    Word.InsertLocation.end\n  );\n  await context.sync();\n}\n\n// --- Change 2 --------------------------------------------------------\n// \"...who have similar questions. Additionally, ...\" ->\n// \"...who have similar questions, as seen in '\u2018. Additionally, ...\"\nlet results2 = body.search(\"who have similar questions\", { matchCase: true, matchWholeWord: false });\nresults2.load(\"items\");\nawait context.sync();\nif (results2.items.length > 0) {\n  results2.items[0].insertText(\n    \", as seen in \\u2019\\u2018\",\n    Word.InsertLocation.end\n  );\n  await context.sync();\n}\n\n// --- Change 3 --------------------------------------------------------\n// \"...to prevent any doubts from lingering.\" ->\n// \"...to prevent any doubts from lingering. I also am proud with how much I\n//  learned over this quarter and how I was able to still have fun with my\n//  assignments especially through making the quarto files have a theme I\n//  liked and incorporating different colors to make the file look more\n//  appealing.\"\n// This sentence is the very end of its paragraph (\"Attention to Personal\n// Goals\"), so insert directly at the paragraph's end.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\nconst goalPara = paragraphs.items.find(\n  (p) => p.text.indexOf(\"to prevent any doubts from lingering.\") !== -1\n);\nif (goalPara) {\n  goalPara.insertText(\n    \" I also am proud with how much I learned over this quarter and how I was able to still have fun with my assignments especially through making the quarto files have a theme I liked and incorporating different colors to make the file look more appealing.\",\n    Word.InsertLocation.end\n  );\n  await context.sync();\n}\n", "ps1": "# Apply the three textual insertions described by the diff using the Word\n# COM object model. Each change appends additional sentence(s)/phrase(s)\n# immediately after an existing, uniquely-identifiable piece of text.\n\n$d = $word.ActiveDocument\n\n# --- Change 1 ----------------------------------------------------------\n# \"...where I had to revise multiple times.\" ->\n# \"...where I had to revise multiple times, and I am proud of the revised\n#  results I was able to get at the end.\"\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"where I had to revise multiple times\"\n$rng.Find.MatchCase = $true\n$found = $rng.Find.Execute()\nif ($found) {\n  $rng.Collapse(0)  # wdCollapseEnd\n  $rng.InsertAfter(\", and I am proud of the revised results I was able to get at the end\")\n}\n\n# --- Change 2 ----------------------------------------------------------\n# \"...who have similar questions. Additionally, ...\" ->\n# \"...who have similar questions, as seen in '`. Additionally, ...\"\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"who have similar questions\"\n$rng2.Find.MatchCase = $true\n$found2 = $rng2.Find.Execute()\nif ($found2) {\n  $rng2.Collapse(0)  # wdCollapseEnd\n  $rng2.InsertAfter(\", as seen in \" + [char]0x2019 + [char]0x2018)\n}\n\n# --- Change 3 ----------------------------------------------------------\n# \"...to prevent any doubts from lingering.\" ->\n# \"...to prevent any doubts from lingering. I also am proud with how much I\n#  learned over this quarter and how I was able to still have fun with my\n#  assignments especially through making the quarto files have a theme I\n#  liked and incorporating different colors to make the file look more\n#  appealing.\"\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Text = \"to prevent any doubts from lingering.\"\n$rng3.Find.MatchCase = $true\n$found3 = $rng3.Find.Execute()\nif ($found3) {\n  $rng3.Collapse(0)  # wdCollapseEnd\n  $rng3.InsertAfter(\" I also am proud with how much I learned over this quarter and how I was able to still have fun with my assignments especially through making the quarto files have a theme I liked and incorporating different colors to make the file look more appealing.\")\n}\n"}
